# The document starts with 3 paragraphs:
#   1. "PBI test documents 1"
#   2. "This document is for testing purpose" + " "
#   3. "This file is modified" + bookmark "_GoBack"
#
# It needs to become 4 paragraphs:
#   1. "File is modified by " + "Varun" + " on 17-Aug-23"
#   2. bookmark "_GoBack" + "PBI test documents 1"   (the old title, now 2nd)
#   3. "This document is for testing purpose" + " "  (unchanged, shifts down)
#   4. "This file is modified" + " " + "This" + " file is modified"
#      (extended text, "_GoBack" bookmark removed from here)

$d = $word.ActiveDocument

# --- 1. Rewrite paragraph 1's text: "PBI test documents 1" -> "File is modified by Varun on 17-Aug-23" ---
$p1 = $d.Paragraphs.Item(1)
$xmlPara1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">File is modified by </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Varun</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> on 17-Aug-23</w:t></w:r>
</w:p>
</w:body>
</w:document>
'@
$p1.Range.InsertXML($xmlPara1)

# --- 2. Insert a brand new paragraph right after paragraph 1 that carries the
#        original title text, with the "_GoBack" bookmark collapsed at its start ---
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item(2)
$xmlPara2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>PBI test documents 1</w:t></w:r>
</w:p>
</w:body>
</w:document>
'@
$p2.Range.InsertXML($xmlPara2)

# --- 3. Last paragraph ("This file is modified ..."): drop its old "_GoBack"
#        bookmark (it moved to paragraph 2) and extend the sentence ---
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$xmlParaLast = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>This file is modified</w:t></w:r>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>This</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> file is modified</w:t></w:r>
</w:p>
</w:body>
</w:document>
'@
$pLast.Range.InsertXML($xmlParaLast)

Write-Host "Final paragraph count:" $d.Paragraphs.Count
foreach ($p in $d.Paragraphs) {
    Write-Host "PARA: [" $p.Range.Text "]"
}
